# Apply cell value updates to match the updated cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these columns are treated as plain text so that values such as
# "1.018" or "0.000009019" are preserved exactly as strings (not reinterpreted
# as numbers / dates / scientific notation) when written back.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.764.54"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.856.69"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "1.018"
$ws.Range("E4").Value = "  -1.94%  "
$ws.Range("D5").Value = "320.72"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").Value = "0.4378"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "0.3783"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "0.07427"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "21.56"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "1.857.83"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "6.785"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "5.490"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "0.07142"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "88.20"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").Value = "1.022"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "0.000009019"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "1.017"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "27.760.80"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "5.273"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "11.15"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "2.096.03"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "2.032"
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("D26").Value = "156.84"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "18.70"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "5.455"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "1.988"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "120.70"
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("D31").Value = "0.09043"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "1.228"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "0.7692"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "3.011"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("D35").Value = "4.557"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "1.137"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "0.01979"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").Value = "0.05305"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "2.865"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").Value = "0.5186"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "6.964"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").Value = "0.1677"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "8.707"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "110.13"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.72"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "0.4735"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "1.019"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "1.851"
$ws.Range("E51").Value = "  -0.63%  "
